$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Update the MSME enterprise percentage figures (text values, matching the
# original cells which stored these numbers as text strings).
$ws.Range("B10").Value = "'88.06"
$ws.Range("B10").Style = "Normal"

$ws.Range("C10").Value = "'11.69"
$ws.Range("C10").Style = "Normal"

$ws.Range("D10").Value = "'99.75"
$ws.Range("D10").Style = "Normal"
